$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 2687.5
$ws.Range("J40").Value = 2916.6667
$ws.Range("L40").Value = 2916.6667
$ws.Range("N40").Value = -3266.6667
$ws.Range("H86").Value = 37504050
$ws.Range("I86").Value = 54549936
$ws.Range("J86").Value = 3100
$ws.Range("K86").Value = 54549936
$ws.Range("L86").Value = 3100
$ws.Range("M86").Value = -54548813
$ws.Range("N86").Value = -5346
$ws.Range("H89").Value = 37504050
$ws.Range("I89").Value = 54549936
$ws.Range("J89").Value = 3100
$ws.Range("K89").Value = 272749680
$ws.Range("L89").Value = 15500
$ws.Range("M89").Value = -272744064
$ws.Range("N89").Value = -26732
$ws.Range("H92").Value = 15873772
$ws.Range("I92").Value = 20833816
$ws.Range("J92").Value = 1630
$ws.Range("K92").Value = 20833816
$ws.Range("L92").Value = 1630
$ws.Range("M92").Value = -20832568
$ws.Range("N92").Value = -4126
$ws.Range("H113").Value = 0
$ws.Range("I113").Value = 0
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 0
$ws.Range("L113").Value = 0
$ws.Range("H132").Value = 4462.643
$ws.Range("I132").Value = 4497.96
$ws.Range("J132").Value = 4168.3335
$ws.Range("K132").Value = 13493.88
$ws.Range("L132").Value = 12505.0005
$ws.Range("M132").Value = -10963.88
$ws.Range("N132").Value = -17565.0005
$ws.Range("H135").Value = 1535.1364
$ws.Range("I135").Value = 762.375
$ws.Range("J135").Value = 3595.8333
$ws.Range("K135").Value = 6861.375
$ws.Range("L135").Value = 32362.4997
$ws.Range("M135").Value = -4326.375
$ws.Range("N135").Value = -37432.4997
$ws.Range("H137").Value = 2815.38
$ws.Range("I137").Value = 2354.9707
$ws.Range("J137").Value = 3793.75
$ws.Range("K137").Value = 7064.9121
$ws.Range("L137").Value = 11381.25
$ws.Range("M137").Value = -4514.9121
$ws.Range("N137").Value = -16481.25
$ws.Range("H138").Value = 2593.3438
$ws.Range("I138").Value = 1582.4517
$ws.Range("J138").Value = 3542.9697
$ws.Range("K138").Value = 4747.355100000001
$ws.Range("L138").Value = 10628.9091
$ws.Range("M138").Value = 392.6448999999993
$ws.Range("N138").Value = -20908.9091
$ws.Range("M113").ClearContents()
$ws.Range("N113").ClearContents()

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1416.2667
$ws.Range("I2").Value = 1543.4445
$ws.Range("J2").Value = 1225.5
$ws.Range("K2").Value = 1543.4445
$ws.Range("L2").Value = 1225.5
$ws.Range("M2").Value = -1430.4445
$ws.Range("N2").Value = -1451.5
$ws.Range("H32").Value = 616927.5600000001
$ws.Range("I32").Value = 667225.6
$ws.Range("J32").Value = 33469.6
$ws.Range("K32").Value = 667225.6
$ws.Range("L32").Value = 33469.6
$ws.Range("M32").Value = -666938.6
$ws.Range("N32").Value = -34043.6
$ws.Range("H61").Value = 3352.3333
$ws.Range("I61").Value = 2775.8667
$ws.Range("K61").Value = 2775.8667
$ws.Range("M61").Value = -2563.8667
$ws.Range("H116").Value = 1416.2667
$ws.Range("I116").Value = 1543.4445
$ws.Range("J116").Value = 1225.5
$ws.Range("K116").Value = 1543.4445
$ws.Range("L116").Value = 1225.5
$ws.Range("M116").Value = 750.5554999999999
$ws.Range("N116").Value = -5813.5
$ws.Range("H122").Value = 59917.883
$ws.Range("I122").Value = 72293.14
$ws.Range("J122").Value = 2166.6667
$ws.Range("K122").Value = 216879.42
$ws.Range("L122").Value = 6500.000100000001
$ws.Range("M122").Value = -214429.42
$ws.Range("N122").Value = -11400.0001
$ws.Range("H136").Value = 3352.3333
$ws.Range("I136").Value = 2775.8667
$ws.Range("K136").Value = 8327.6001
$ws.Range("M136").Value = -5777.6001
$ws.Range("H138").Value = 47429
$ws.Range("J138").Value = 47429
$ws.Range("L138").Value = 47429
$ws.Range("N138").Value = -57709

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1416.2667
$ws.Range("I3").Value = 1543.4445
$ws.Range("J3").Value = 1225.5
$ws.Range("K3").Value = 1543.4445
$ws.Range("L3").Value = 1225.5
$ws.Range("M3").Value = -1429.4445
$ws.Range("N3").Value = -1453.5
$ws.Range("H94").Value = 2587.1428
$ws.Range("I94").Value = 1600
$ws.Range("K94").Value = 1600
$ws.Range("M94").Value = -1149
$ws.Range("H99").Value = 2416.6667
$ws.Range("I99").Value = 1800
$ws.Range("J99").Value = 2725
$ws.Range("K99").Value = 1800
$ws.Range("L99").Value = 2725
$ws.Range("M99").Value = -302
$ws.Range("N99").Value = -5721

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H3").Value = 4280
$ws.Range("J3").Value = 4280
$ws.Range("L3").Value = 4280
$ws.Range("N3").Value = -4506
$ws.Range("H31").Value = 6785.125
$ws.Range("I31").Value = 1275.5652
$ws.Range("J31").Value = 11853.92
$ws.Range("K31").Value = 1275.5652
$ws.Range("L31").Value = 11853.92
$ws.Range("M31").Value = -980.5652
$ws.Range("N31").Value = -12443.92
$ws.Range("H34").Value = 6785.125
$ws.Range("I34").Value = 1275.5652
$ws.Range("J34").Value = 11853.92
$ws.Range("K34").Value = 1275.5652
$ws.Range("L34").Value = 11853.92
$ws.Range("M34").Value = -1073.5652
$ws.Range("N34").Value = -12257.92
$ws.Range("H58").Value = 1583.8823
$ws.Range("I58").Value = 1328
$ws.Range("J58").Value = 1949.4286
$ws.Range("K58").Value = 1328
$ws.Range("L58").Value = 1949.4286
$ws.Range("M58").Value = -1125
$ws.Range("N58").Value = -2355.4286
$ws.Range("H100").Value = 40000
$ws.Range("J100").Value = 40000
$ws.Range("L100").Value = 40000
$ws.Range("N100").Value = -42164
$ws.Range("H105").Value = 400
$ws.Range("I105").Value = 400
$ws.Range("K105").Value = 400
$ws.Range("M105").Value = 1347
$ws.Range("H134").Value = 5159.9653
$ws.Range("I134").Value = 5377.1304
$ws.Range("J134").Value = 4327.5
$ws.Range("K134").Value = 16131.3912
$ws.Range("L134").Value = 12982.5
$ws.Range("M134").Value = -13596.3912
$ws.Range("N134").Value = -18052.5
$ws.Range("H136").Value = 1583.8823
$ws.Range("I136").Value = 1328
$ws.Range("J136").Value = 1949.4286
$ws.Range("K136").Value = 3984
$ws.Range("L136").Value = 5848.2858
$ws.Range("M136").Value = -1434
$ws.Range("N136").Value = -10948.2858

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 708.85187
$ws.Range("J5").Value = 825.3333
$ws.Range("L5").Value = 2475.9999
$ws.Range("N5").Value = -2699.9999
$ws.Range("H31").Value = 1756.1428
$ws.Range("J31").Value = 1660.4615
$ws.Range("L31").Value = 4981.3845
$ws.Range("N31").Value = -5557.3845
$ws.Range("H111").Value = 5928.4287
$ws.Range("I111").Value = 624.75
$ws.Range("J111").Value = 13000
$ws.Range("K111").Value = 1874.25
$ws.Range("L111").Value = 39000
$ws.Range("M111").Value = 1192.75
$ws.Range("N111").Value = -45134
$ws.Range("H131").Value = 1016.2105
$ws.Range("I131").Value = 654.8333
$ws.Range("J131").Value = 1183
$ws.Range("K131").Value = 1964.4999
$ws.Range("L131").Value = 3549
$ws.Range("M131").Value = 3075.5001
$ws.Range("N131").Value = -13629
$ws.Range("H135").Value = 708.85187
$ws.Range("J135").Value = 825.3333
$ws.Range("L135").Value = 7427.9997
$ws.Range("N135").Value = -12497.9997

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 1743.7778
$ws.Range("I122").Value = 1338.8
$ws.Range("K122").Value = 4016.4
$ws.Range("M122").Value = -1566.4

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 7248066
$ws.Range("I136").Value = 1622.6154
$ws.Range("J136").Value = 16668442
$ws.Range("K136").Value = 4867.8462
$ws.Range("L136").Value = 50005326
$ws.Range("M136").Value = -2317.8462
$ws.Range("N136").Value = -50010426

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 3213.5
$ws.Range("I136").Value = 2748.8948
$ws.Range("K136").Value = 8246.6844
$ws.Range("M136").Value = -5696.6844
